$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "001" to "002" (text, keep no special number formatting).
# Writing a leading-zero numeric-looking string normally gets coerced to a
# number by the engine, so force it in as text (leading apostrophe) and then
# reset the style back to Normal so no stray style index is left on the cell.
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# K2 stays "001" - no change needed.

# N2: report date updated.
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Updated absolute figures.
$ws.Range("O2").Value = 465486349.12
$ws.Range("P2").Value = 140201801.52
$ws.Range("Q2").Value = 12302265.4
$ws.Range("S2").Value = 145693570.47
$ws.Range("U2").Value = 81077923.83
$ws.Range("W2").Value = 147109632.17
$ws.Range("X2").Value = 90065070.09999999
$ws.Range("Z2").Value = 122890.12
$ws.Range("AB2").Value = 318376716.95
$ws.Range("AF2").Value = 213.7804359889
$ws.Range("AG2").Value = 31.6034256317

# Ratio columns are cleared (no longer populated in the new snapshot).
$ws.Range("R2").ClearContents()
$ws.Range("T2").ClearContents()
$ws.Range("V2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").ClearContents()
